$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-31 with new TPM-derived values (columns E-T only; A-D unchanged)
# Row 2
$ws.Range("G2").Value = 1.442444333333333
$ws.Range("H2").Value = 4.327332999999999
$ws.Range("I2").Value = 0.1935924704900452
$ws.Range("J2").Value = 0.1935924704900451
$ws.Range("M2").Value = 0.978762
$ws.Range("N2").Value = 2.936286
$ws.Range("O2").Value = 0.3819465121442868
$ws.Range("P2").Value = 0.3819465121442868
$ws.Range("Q2").Value = 1.411809700582
$ws.Range("R2").Value = 12.706287305238
$ws.Range("S2").Value = 0.07394196888106853
$ws.Range("T2").Value = 0.0739419688810685
# Row 3
$ws.Range("G3").Value = 1.442444333333333
$ws.Range("H3").Value = 4.327332999999999
$ws.Range("I3").Value = 0.1935924704900452
$ws.Range("J3").Value = 0.1935924704900451
$ws.Range("O3").Value = 0.0406458950147437
$ws.Range("P3").Value = 0.04064589501474371
$ws.Range("Q3").Value = 0.1502416360565555
$ws.Range("R3").Value = 1.352174724509
$ws.Range("S3").Value = 0.007868739231183244
$ws.Range("T3").Value = 0.007868739231183244
# Row 4
$ws.Range("G4").Value = 1.442444333333333
$ws.Range("H4").Value = 4.327332999999999
$ws.Range("I4").Value = 0.1935924704900452
$ws.Range("J4").Value = 0.1935924704900451
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3139526666666667
$ws.Range("N4").Value = 0.9418580000000001
$ws.Range("O4").Value = 0.1225151017425393
$ws.Range("P4").Value = 0.1225151017425393
$ws.Range("Q4").Value = 0.4528592449682223
$ws.Range("R4").Value = 4.075733204714
$ws.Range("S4").Value = 0.02371800121867742
$ws.Range("T4").Value = 0.02371800121867742
# Row 5
$ws.Range("G5").Value = 1.442444333333333
$ws.Range("H5").Value = 4.327332999999999
$ws.Range("I5").Value = 0.1935924704900452
$ws.Range("J5").Value = 0.1935924704900451
$ws.Range("M5").Value = 0.5829876666666666
$ws.Range("N5").Value = 1.748963
$ws.Range("O5").Value = 0.2275017888991087
$ws.Range("P5").Value = 0.2275017888991087
$ws.Range("Q5").Value = 0.8409272561865554
$ws.Range("R5").Value = 7.568345305678998
$ws.Range("S5").Value = 0.04404263335388318
$ws.Range("T5").Value = 0.04404263335388318
# Row 6
$ws.Range("G6").Value = 1.442444333333333
$ws.Range("H6").Value = 4.327332999999999
$ws.Range("I6").Value = 0.1935924704900452
$ws.Range("J6").Value = 0.1935924704900451
$ws.Range("M6").Value = 0.3019996666666667
$ws.Range("N6").Value = 0.905999
$ws.Range("O6").Value = 0.117850631054404
$ws.Range("P6").Value = 0.117850631054404
$ws.Range("Q6").Value = 0.4356177078518889
$ws.Range("R6").Value = 3.920559370667
$ws.Range("S6").Value = 0.02281499481463291
$ws.Range("T6").Value = 0.0228149948146329
# Row 7
$ws.Range("G7").Value = 1.442444333333333
$ws.Range("H7").Value = 4.327332999999999
$ws.Range("I7").Value = 0.1935924704900452
$ws.Range("J7").Value = 0.1935924704900451
$ws.Range("M7").Value = 0.2807033333333333
$ws.Range("N7").Value = 0.8421099999999999
$ws.Range("O7").Value = 0.1095400711449175
$ws.Range("P7").Value = 0.1095400711449176
$ws.Range("Q7").Value = 0.4048989325144444
$ws.Range("R7").Value = 3.644090392629999
$ws.Range("S7").Value = 0.0212061329905999
$ws.Range("T7").Value = 0.0212061329905999
# Row 8
$ws.Range("I8").Value = 0.3798216115782563
$ws.Range("J8").Value = 0.3798216115782563
$ws.Range("M8").Value = 0.978762
$ws.Range("N8").Value = 2.936286
$ws.Range("O8").Value = 0.3819465121442868
$ws.Range("P8").Value = 0.3819465121442868
$ws.Range("Q8").Value = 2.76992092905
$ws.Range("R8").Value = 24.92928836145
$ws.Range("S8").Value = 0.1450715397793371
$ws.Range("T8").Value = 0.1450715397793371
# Row 9
$ws.Range("I9").Value = 0.3798216115782563
$ws.Range("J9").Value = 0.3798216115782563
$ws.Range("O9").Value = 0.0406458950147437
$ws.Range("P9").Value = 0.04064589501474371
$ws.Range("S9").Value = 0.01543818934854057
$ws.Range("T9").Value = 0.01543818934854057
# Row 10
$ws.Range("I10").Value = 0.3798216115782563
$ws.Range("J10").Value = 0.3798216115782563
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3139526666666667
$ws.Range("N10").Value = 0.9418580000000001
$ws.Range("O10").Value = 0.1225151017425393
$ws.Range("P10").Value = 0.1225151017425393
$ws.Range("Q10").Value = 0.8884938954833336
$ws.Range("R10").Value = 7.996445059350002
$ws.Range("S10").Value = 0.04653388338652532
$ws.Range("T10").Value = 0.04653388338652531
# Row 11
$ws.Range("I11").Value = 0.3798216115782563
$ws.Range("J11").Value = 0.3798216115782563
$ws.Range("M11").Value = 0.5829876666666666
$ws.Range("N11").Value = 1.748963
$ws.Range("O11").Value = 0.2275017888991087
$ws.Range("P11").Value = 0.2275017888991087
$ws.Range("Q11").Value = 1.649869671358333
$ws.Range("R11").Value = 14.848827042225
$ws.Range("S11").Value = 0.08641009609659572
$ws.Range("T11").Value = 0.08641009609659572
# Row 12
$ws.Range("I12").Value = 0.3798216115782563
$ws.Range("J12").Value = 0.3798216115782563
$ws.Range("M12").Value = 0.3019996666666667
$ws.Range("N12").Value = 0.905999
$ws.Range("O12").Value = 0.117850631054404
$ws.Range("P12").Value = 0.117850631054404
$ws.Range("Q12").Value = 0.8546666066583335
$ws.Range("R12").Value = 7.691999459925001
$ws.Range("S12").Value = 0.04476221661259823
$ws.Range("T12").Value = 0.04476221661259822
# Row 13
$ws.Range("I13").Value = 0.3798216115782563
$ws.Range("J13").Value = 0.3798216115782563
$ws.Range("M13").Value = 0.2807033333333333
$ws.Range("N13").Value = 0.8421099999999999
$ws.Range("O13").Value = 0.1095400711449175
$ws.Range("P13").Value = 0.1095400711449176
$ws.Range("Q13").Value = 0.7943974509166667
$ws.Range("R13").Value = 7.14957705825
$ws.Range("S13").Value = 0.04160568635465944
$ws.Range("T13").Value = 0.04160568635465944
# Row 14
$ws.Range("I14").Value = 0.05301717064836096
$ws.Range("J14").Value = 0.05301717064836094
$ws.Range("M14").Value = 0.978762
$ws.Range("N14").Value = 2.936286
$ws.Range("O14").Value = 0.3819465121442868
$ws.Range("P14").Value = 0.3819465121442868
$ws.Range("Q14").Value = 0.386637742828
$ws.Range("R14").Value = 3.479739685452
$ws.Range("S14").Value = 0.02024972341289993
$ws.Range("T14").Value = 0.02024972341289992
# Row 15
$ws.Range("I15").Value = 0.05301717064836096
$ws.Range("J15").Value = 0.05301717064836094
$ws.Range("O15").Value = 0.0406458950147437
$ws.Range("P15").Value = 0.04064589501474371
$ws.Range("S15").Value = 0.002154930352152031
$ws.Range("T15").Value = 0.00215493035215203
# Row 16
$ws.Range("I16").Value = 0.05301717064836096
$ws.Range("J16").Value = 0.05301717064836094
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.3139526666666667
$ws.Range("N16").Value = 0.9418580000000001
$ws.Range("O16").Value = 0.1225151017425393
$ws.Range("P16").Value = 0.1225151017425393
$ws.Range("Q16").Value = 0.1240198847062222
$ws.Range("R16").Value = 1.116178962356
$ws.Range("S16").Value = 0.006495404056085511
$ws.Range("T16").Value = 0.006495404056085509
# Row 17
$ws.Range("I17").Value = 0.05301717064836096
$ws.Range("J17").Value = 0.05301717064836094
$ws.Range("M17").Value = 0.5829876666666666
$ws.Range("N17").Value = 1.748963
$ws.Range("O17").Value = 0.2275017888991087
$ws.Range("P17").Value = 0.2275017888991087
$ws.Range("Q17").Value = 0.2302960633295555
$ws.Range("R17").Value = 2.072664569966
$ws.Range("S17").Value = 0.01206150116487143
$ws.Range("T17").Value = 0.01206150116487143
# Row 18
$ws.Range("I18").Value = 0.05301717064836096
$ws.Range("J18").Value = 0.05301717064836094
$ws.Range("M18").Value = 0.3019996666666667
$ws.Range("N18").Value = 0.905999
$ws.Range("O18").Value = 0.117850631054404
$ws.Range("P18").Value = 0.117850631054404
$ws.Range("Q18").Value = 0.1192981229908889
$ws.Range("R18").Value = 1.073683106918
$ws.Range("S18").Value = 0.006248107017628363
$ws.Range("T18").Value = 0.006248107017628362
# Row 19
$ws.Range("I19").Value = 0.05301717064836096
$ws.Range("J19").Value = 0.05301717064836094
$ws.Range("M19").Value = 0.2807033333333333
$ws.Range("N19").Value = 0.8421099999999999
$ws.Range("O19").Value = 0.1095400711449175
$ws.Range("P19").Value = 0.1095400711449176
$ws.Range("Q19").Value = 0.1108854892244444
$ws.Range("R19").Value = 0.9979694030199999
$ws.Range("S19").Value = 0.005807504644723693
$ws.Range("T19").Value = 0.005807504644723693
# Row 20
$ws.Range("G20").Value = 2.258134333333333
$ws.Range("H20").Value = 6.774403
$ws.Range("I20").Value = 0.3030673657112068
$ws.Range("J20").Value = 0.3030673657112067
$ws.Range("M20").Value = 0.978762
$ws.Range("N20").Value = 2.936286
$ws.Range("O20").Value = 0.3819465121442868
$ws.Range("P20").Value = 0.3819465121442868
$ws.Range("Q20").Value = 2.210176076362
$ws.Range("R20").Value = 19.891584687258
$ws.Range("S20").Value = 0.1157555232781525
$ws.Range("T20").Value = 0.1157555232781524
# Row 21
$ws.Range("G21").Value = 2.258134333333333
$ws.Range("H21").Value = 6.774403
$ws.Range("I21").Value = 0.3030673657112068
$ws.Range("J21").Value = 0.3030673657112067
$ws.Range("O21").Value = 0.0406458950147437
$ws.Range("P21").Value = 0.04064589501474371
$ws.Range("Q21").Value = 0.2352020031798888
$ws.Range("R21").Value = 2.116818028619
$ws.Range("S21").Value = 0.01231844432909265
$ws.Range("T21").Value = 0.01231844432909265
# Row 22
$ws.Range("G22").Value = 2.258134333333333
$ws.Range("H22").Value = 6.774403
$ws.Range("I22").Value = 0.3030673657112068
$ws.Range("J22").Value = 0.3030673657112067
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.3139526666666667
$ws.Range("N22").Value = 0.9418580000000001
$ws.Range("O22").Value = 0.1225151017425393
$ws.Range("P22").Value = 0.1225151017425393
$ws.Range("Q22").Value = 0.7089472956415556
$ws.Range("R22").Value = 6.380525660774
$ws.Range("S22").Value = 0.03713032914495187
$ws.Range("T22").Value = 0.03713032914495186
# Row 23
$ws.Range("G23").Value = 2.258134333333333
$ws.Range("H23").Value = 6.774403
$ws.Range("I23").Value = 0.3030673657112068
$ws.Range("J23").Value = 0.3030673657112067
$ws.Range("M23").Value = 0.5829876666666666
$ws.Range("N23").Value = 1.748963
$ws.Range("O23").Value = 0.2275017888991087
$ws.Range("P23").Value = 0.2275017888991087
$ws.Range("Q23").Value = 1.316464466009889
$ws.Range("R23").Value = 11.848180194089
$ws.Range("S23").Value = 0.06894836785623994
$ws.Range("T23").Value = 0.06894836785623994
# Row 24
$ws.Range("G24").Value = 2.258134333333333
$ws.Range("H24").Value = 6.774403
$ws.Range("I24").Value = 0.3030673657112068
$ws.Range("J24").Value = 0.3030673657112067
$ws.Range("M24").Value = 0.3019996666666667
$ws.Range("N24").Value = 0.905999
$ws.Range("O24").Value = 0.117850631054404
$ws.Range("P24").Value = 0.117850631054404
$ws.Range("Q24").Value = 0.6819558159552221
$ws.Range("R24").Value = 6.137602343597
$ws.Range("S24").Value = 0.03571668030106156
$ws.Range("T24").Value = 0.03571668030106156
# Row 25
$ws.Range("G25").Value = 2.258134333333333
$ws.Range("H25").Value = 6.774403
$ws.Range("I25").Value = 0.3030673657112068
$ws.Range("J25").Value = 0.3030673657112067
$ws.Range("M25").Value = 0.2807033333333333
$ws.Range("N25").Value = 0.8421099999999999
$ws.Range("O25").Value = 0.1095400711449175
$ws.Range("P25").Value = 0.1095400711449176
$ws.Range("Q25").Value = 0.633865834481111
$ws.Range("R25").Value = 5.704792510329999
$ws.Range("S25").Value = 0.03319802080170833
$ws.Range("T25").Value = 0.03319802080170833
# Row 26
$ws.Range("G26").Value = 0.337485
$ws.Range("H26").Value = 1.012455
$ws.Range("I26").Value = 0.04529433364846171
$ws.Range("J26").Value = 0.0452943336484617
$ws.Range("M26").Value = 0.978762
$ws.Range("N26").Value = 2.936286
$ws.Range("O26").Value = 0.3819465121442868
$ws.Range("P26").Value = 0.3819465121442868
$ws.Range("Q26").Value = 0.3303174935700001
$ws.Range("R26").Value = 2.97285744213
$ws.Range("S26").Value = 0.01730001275692956
$ws.Range("T26").Value = 0.01730001275692956
# Row 27
$ws.Range("G27").Value = 0.337485
$ws.Range("H27").Value = 1.012455
$ws.Range("I27").Value = 0.04529433364846171
$ws.Range("J27").Value = 0.0452943336484617
$ws.Range("O27").Value = 0.0406458950147437
$ws.Range("P27").Value = 0.04064589501474371
$ws.Range("Q27").Value = 0.035151650135
$ws.Range("R27").Value = 0.316364851215
$ws.Range("S27").Value = 0.001841028730238147
$ws.Range("T27").Value = 0.001841028730238147
# Row 28
$ws.Range("G28").Value = 0.337485
$ws.Range("H28").Value = 1.012455
$ws.Range("I28").Value = 0.04529433364846171
$ws.Range("J28").Value = 0.0452943336484617
$ws.Range("K28").Value = 2
$ws.Range("L28").Value = 0.6666666666666666
$ws.Range("M28").Value = 0.3139526666666667
$ws.Range("N28").Value = 0.9418580000000001
$ws.Range("O28").Value = 0.1225151017425393
$ws.Range("P28").Value = 0.1225151017425393
$ws.Range("Q28").Value = 0.10595431571
$ws.Range("R28").Value = 0.9535888413900002
$ws.Range("S28").Value = 0.005549239895301807
$ws.Range("T28").Value = 0.005549239895301806
# Row 29
$ws.Range("G29").Value = 0.337485
$ws.Range("H29").Value = 1.012455
$ws.Range("I29").Value = 0.04529433364846171
$ws.Range("J29").Value = 0.0452943336484617
$ws.Range("M29").Value = 0.5829876666666666
$ws.Range("N29").Value = 1.748963
$ws.Range("O29").Value = 0.2275017888991087
$ws.Range("P29").Value = 0.2275017888991087
$ws.Range("Q29").Value = 0.196749592685
$ws.Range("R29").Value = 1.770746334165
$ws.Range("S29").Value = 0.01030454193201813
$ws.Range("T29").Value = 0.01030454193201813
# Row 30
$ws.Range("G30").Value = 0.337485
$ws.Range("H30").Value = 1.012455
$ws.Range("I30").Value = 0.04529433364846171
$ws.Range("J30").Value = 0.0452943336484617
$ws.Range("M30").Value = 0.3019996666666667
$ws.Range("N30").Value = 0.905999
$ws.Range("O30").Value = 0.117850631054404
$ws.Range("P30").Value = 0.117850631054404
$ws.Range("Q30").Value = 0.101920357505
$ws.Range("R30").Value = 0.9172832175450001
$ws.Range("S30").Value = 0.005337965803659937
$ws.Range("T30").Value = 0.005337965803659936
# Row 31
$ws.Range("G31").Value = 0.337485
$ws.Range("H31").Value = 1.012455
$ws.Range("I31").Value = 0.04529433364846171
$ws.Range("J31").Value = 0.0452943336484617
$ws.Range("M31").Value = 0.2807033333333333
$ws.Range("N31").Value = 0.8421099999999999
$ws.Range("O31").Value = 0.1095400711449175
$ws.Range("P31").Value = 0.1095400711449176
$ws.Range("Q31").Value = 0.09473316445
$ws.Range("R31").Value = 0.85259848005
$ws.Range("S31").Value = 0.004961544530314128
$ws.Range("T31").Value = 0.004961544530314128

# Add new rows 32-37 for Resolving-Mac as sending cluster
# Row 32
$ws.Range("A32").Value = "Resolving-Mac"
$ws.Range("B32").Value = "Gdf11"
$ws.Range("C32").Value = "Acvr2b"
$ws.Range("D32").Value = "ECs"
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = 0.6666666666666666
$ws.Range("G32").Value = 0.187816
$ws.Range("H32").Value = 0.563448
$ws.Range("I32").Value = 0.02520704792366915
$ws.Range("J32").Value = 0.02520704792366914
$ws.Range("K32").Value = 3
$ws.Range("L32").Value = 1
$ws.Range("M32").Value = 0.978762
$ws.Range("N32").Value = 2.936286
$ws.Range("O32").Value = 0.3819465121442868
$ws.Range("P32").Value = 0.3819465121442868
$ws.Range("Q32").Value = 0.183827163792
$ws.Range("R32").Value = 1.654444474128
$ws.Range("S32").Value = 0.009627744035899319
$ws.Range("T32").Value = 0.009627744035899315
# Row 33
$ws.Range("A33").Value = "Resolving-Mac"
$ws.Range("B33").Value = "Gdf11"
$ws.Range("C33").Value = "Acvr2b"
$ws.Range("D33").Value = "FAPs"
$ws.Range("E33").Value = 2
$ws.Range("F33").Value = 0.6666666666666666
$ws.Range("G33").Value = 0.187816
$ws.Range("H33").Value = 0.563448
$ws.Range("I33").Value = 0.02520704792366915
$ws.Range("J33").Value = 0.02520704792366914
$ws.Range("K33").Value = 2
$ws.Range("L33").Value = 0.6666666666666666
$ws.Range("M33").Value = 0.1041576666666667
$ws.Range("N33").Value = 0.312473
$ws.Range("O33").Value = 0.0406458950147437
$ws.Range("P33").Value = 0.04064589501474371
$ws.Range("Q33").Value = 0.01956247632266666
$ws.Range("R33").Value = 0.176062286904
$ws.Range("S33").Value = 0.001024563023537069
$ws.Range("T33").Value = 0.001024563023537069
# Row 34
$ws.Range("A34").Value = "Resolving-Mac"
$ws.Range("B34").Value = "Gdf11"
$ws.Range("C34").Value = "Acvr2b"
$ws.Range("D34").Value = "Inflammatory-Mac"
$ws.Range("E34").Value = 2
$ws.Range("F34").Value = 0.6666666666666666
$ws.Range("G34").Value = 0.187816
$ws.Range("H34").Value = 0.563448
$ws.Range("I34").Value = 0.02520704792366915
$ws.Range("J34").Value = 0.02520704792366914
$ws.Range("K34").Value = 2
$ws.Range("L34").Value = 0.6666666666666666
$ws.Range("M34").Value = 0.3139526666666667
$ws.Range("N34").Value = 0.9418580000000001
$ws.Range("O34").Value = 0.1225151017425393
$ws.Range("P34").Value = 0.1225151017425393
$ws.Range("Q34").Value = 0.05896533404266667
$ws.Range("R34").Value = 0.530688006384
$ws.Range("S34").Value = 0.003088244040997389
$ws.Range("T34").Value = 0.003088244040997389
# Row 35
$ws.Range("A35").Value = "Resolving-Mac"
$ws.Range("B35").Value = "Gdf11"
$ws.Range("C35").Value = "Acvr2b"
$ws.Range("D35").Value = "MuSCs"
$ws.Range("E35").Value = 2
$ws.Range("F35").Value = 0.6666666666666666
$ws.Range("G35").Value = 0.187816
$ws.Range("H35").Value = 0.563448
$ws.Range("I35").Value = 0.02520704792366915
$ws.Range("J35").Value = 0.02520704792366914
$ws.Range("K35").Value = 3
$ws.Range("L35").Value = 1
$ws.Range("M35").Value = 0.5829876666666666
$ws.Range("N35").Value = 1.748963
$ws.Range("O35").Value = 0.2275017888991087
$ws.Range("P35").Value = 0.2275017888991087
$ws.Range("Q35").Value = 0.1094944116026666
$ws.Range("R35").Value = 0.9854497044239998
$ws.Range("S35").Value = 0.005734648495500294
$ws.Range("T35").Value = 0.005734648495500294
# Row 36
$ws.Range("A36").Value = "Resolving-Mac"
$ws.Range("B36").Value = "Gdf11"
$ws.Range("C36").Value = "Acvr2b"
$ws.Range("D36").Value = "Neutrophils"
$ws.Range("E36").Value = 2
$ws.Range("F36").Value = 0.6666666666666666
$ws.Range("G36").Value = 0.187816
$ws.Range("H36").Value = 0.563448
$ws.Range("I36").Value = 0.02520704792366915
$ws.Range("J36").Value = 0.02520704792366914
$ws.Range("K36").Value = 2
$ws.Range("L36").Value = 0.6666666666666666
$ws.Range("M36").Value = 0.3019996666666667
$ws.Range("N36").Value = 0.905999
$ws.Range("O36").Value = 0.117850631054404
$ws.Range("P36").Value = 0.117850631054404
$ws.Range("Q36").Value = 0.05672036939466666
$ws.Range("R36").Value = 0.510483324552
$ws.Range("S36").Value = 0.002970666504823013
$ws.Range("T36").Value = 0.002970666504823013
# Row 37
$ws.Range("A37").Value = "Resolving-Mac"
$ws.Range("B37").Value = "Gdf11"
$ws.Range("C37").Value = "Acvr2b"
$ws.Range("D37").Value = "Resolving-Mac"
$ws.Range("E37").Value = 2
$ws.Range("F37").Value = 0.6666666666666666
$ws.Range("G37").Value = 0.187816
$ws.Range("H37").Value = 0.563448
$ws.Range("I37").Value = 0.02520704792366915
$ws.Range("J37").Value = 0.02520704792366914
$ws.Range("K37").Value = 3
$ws.Range("L37").Value = 1
$ws.Range("M37").Value = 0.2807033333333333
$ws.Range("N37").Value = 0.8421099999999999
$ws.Range("O37").Value = 0.1095400711449175
$ws.Range("P37").Value = 0.1095400711449176
$ws.Range("Q37").Value = 0.05272057725333332
$ws.Range("R37").Value = 0.4744851952799999
$ws.Range("S37").Value = 0.002761181822912065
$ws.Range("T37").Value = 0.002761181822912064